$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 value tweaks
$ws.Range("B2").Value = 11.017066881382448
$ws.Range("C2").Value = -1.6054758133320206
$ws.Range("D2").Value = -0.25453656396425117
$ws.Range("E2").Value = -2.5754209949959659

# Row 3 value tweaks
$ws.Range("B3").Value = 3.4386750814914819
$ws.Range("C3").Value = 3.635010254214933
$ws.Range("D3").Value = -0.1021508699954552
$ws.Range("E3").Value = -5.8337590266512791

# Update selection to match new narrower highlighted range
$ws.Range("B1:E3").Select()
